# Dr Khadeejah Project Edited by Danjuma
#
# The heading run "Applied Question: " (with trailing colon + space) is
# changed to "Applied Question" (no trailing colon/space) and the
# paragraph mark that used to terminate that run is pushed down to a new,
# otherwise-empty paragraph inserted right after it. The new paragraph
# inherits the same paragraph formatting (bold, position -30, Times New
# Roman 12pt) that the heading paragraph already carried in its own
# paragraph mark run properties, so a single Find/Replace that swaps the
# trailing ": " for a paragraph break ("^p") reproduces both changes in
# one step: it trims the text AND splits the paragraph, with the new
# empty paragraph naturally inheriting the original pPr/rPr.

$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "Applied Question: ",  # FindText
    $true,                  # MatchCase
    $false,                 # MatchWholeWord
    $false,                 # MatchWildcards
    $false,                 # MatchSoundsLike
    $false,                 # MatchAllWordForms
    $true,                  # Forward
    1,                      # Wrap (wdFindContinue)
    $false,                 # Format
    "Applied Question^p",   # ReplaceWith  (^p = new paragraph mark)
    2                        # Replace (wdReplaceAll)
)

Write-Output ("Applied Question heading split: " + $found)
